# Apply updated cryptocurrency price/volume data as described in the commit diff.
# Note: some "Price" values are plain decimal numbers (e.g. "569.11"); a leading
# apostrophe is used for those so Excel stores them as literal text (matching the
# source data) instead of auto-converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '60.407.56'
$ws.Cells.Item(2, 5).Value = '  +3.30%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.661.16'
$ws.Cells.Item(3, 5).Value = '  +1.35%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  -0.11%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''569.11'
$ws.Cells.Item(5, 5).Value = '  +6.49%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''146.01'
$ws.Cells.Item(6, 5).Value = '  +2.27%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.608'
$ws.Cells.Item(8, 5).Value = '  +7.15%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -2.28%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.104'
$ws.Cells.Item(10, 5).Value = '  +3.68%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +2.42%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '3.106.02'
$ws.Cells.Item(13, 5).Value = '  +0.38%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '60.364.42'
$ws.Cells.Item(14, 5).Value = '  +3.31%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''21.73'
$ws.Cells.Item(15, 5).Value = '  +4.61%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.654.74'
$ws.Cells.Item(16, 5).Value = '  +1.22%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +3.36%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''4.55'
$ws.Cells.Item(18, 5).Value = '  +3.78%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''344.76'
$ws.Cells.Item(19, 5).Value = '  +3.03%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''10.39'
$ws.Cells.Item(20, 5).Value = '  +2.48%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''6.37'
$ws.Cells.Item(21, 5).Value = '  +2.27%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''5.81'
$ws.Cells.Item(22, 5).Value = '  +0.86%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''0.998'
$ws.Cells.Item(23, 5).Value = '  -0.01%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''66.73'
$ws.Cells.Item(24, 5).Value = '  +1.06%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +6.20%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +1.87%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''0.994'
$ws.Cells.Item(27, 5).Value = '  -0.94%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''7.34'
$ws.Cells.Item(28, 5).Value = '  +3.02%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '0.0₃0772'
$ws.Cells.Item(29, 5).Value = '  +4.94%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''1.71'
$ws.Cells.Item(31, 5).Value = '  +4.49%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''6.12'
$ws.Cells.Item(32, 5).Value = '  +4.29%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''156.26'
$ws.Cells.Item(33, 5).Value = '  +3.80%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''4.09'
$ws.Cells.Item(35, 5).Value = '  +5.12%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'SuiNetwork'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(36, 4).Value = '''0.907'
$ws.Cells.Item(36, 5).Value = '  +6.93%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Fetch.AI'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(37, 4).Value = '''0.911'
$ws.Cells.Item(37, 5).Value = '  +12.34%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +6.20%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''37.58'
$ws.Cells.Item(39, 5).Value = '  +1.10%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +6.97%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''308.04'
$ws.Cells.Item(41, 5).Value = '  +9.63%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''3.66'
$ws.Cells.Item(42, 5).Value = '  +2.66%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.996'
$ws.Cells.Item(43, 5).Value = '  -0.30%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Mantle'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(44, 4).Value = '''0.606'
$ws.Cells.Item(44, 5).Value = '  +1.19%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Stellar'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(45, 4).Value = '''0.0980'
$ws.Cells.Item(45, 5).Value = '  +4.79%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +3.04%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''19.37'
$ws.Cells.Item(47, 5).Value = '  +2.00%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''10.66'
$ws.Cells.Item(48, 5).Value = '  -0.36%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''126.08'
$ws.Cells.Item(49, 5).Value = '  +11.98%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''0.0235'
$ws.Cells.Item(50, 5).Value = '  +4.63%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Maker'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(51, 4).Value = '1.969.15'
$ws.Cells.Item(51, 5).Value = '  +1.15%  '
